# Updates weekly price data for Hortaliza - Vega Monumental Concepcion - Alcachofa
# Applies the refreshed dataset values (date, variety, volume, price range,
# weighted average price, commercialization unit, origin, $/Kg price and
# Kg/units columns) for rows 4-36 on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(4, 4).Value = 44496
$ws.Cells.Item(4, 8).Value = 'Madrigal'
$ws.Cells.Item(4, 10).Value = 350
$ws.Cells.Item(4, 11).Value = 7000
$ws.Cells.Item(4, 12).Value = 7500
$ws.Cells.Item(4, 13).Value = 7214
$ws.Cells.Item(4, 14).Value = '$/caja 40 unidades'
$ws.Cells.Item(4, 15).Value = 'Región de Coquimbo'
$ws.Cells.Item(4, 16).Value = 180
$ws.Cells.Item(4, 17).Value = 40
$ws.Cells.Item(5, 4).Value = 44364
$ws.Cells.Item(5, 8).Value = 'Argentina(o)'
$ws.Cells.Item(5, 10).Value = 100
$ws.Cells.Item(5, 11).Value = 19000
$ws.Cells.Item(5, 12).Value = 20000
$ws.Cells.Item(5, 13).Value = 19500
$ws.Cells.Item(5, 14).Value = '$/caja 50 unidades'
$ws.Cells.Item(5, 16).Value = 390
$ws.Cells.Item(5, 17).Value = 50
$ws.Cells.Item(6, 4).Value = 44364
$ws.Cells.Item(6, 8).Value = 'Española'
$ws.Cells.Item(6, 10).Value = 100
$ws.Cells.Item(6, 11).Value = 19000
$ws.Cells.Item(6, 12).Value = 20000
$ws.Cells.Item(6, 13).Value = 19500
$ws.Cells.Item(6, 14).Value = '$/caja 30 unidades'
$ws.Cells.Item(6, 16).Value = 650
$ws.Cells.Item(6, 17).Value = 30
$ws.Cells.Item(7, 4).Value = 44376
$ws.Cells.Item(7, 8).Value = 'Española'
$ws.Cells.Item(7, 10).Value = 100
$ws.Cells.Item(7, 11).Value = 19000
$ws.Cells.Item(7, 12).Value = 20000
$ws.Cells.Item(7, 13).Value = 19500
$ws.Cells.Item(7, 14).Value = '$/caja 30 unidades'
$ws.Cells.Item(7, 16).Value = 650
$ws.Cells.Item(7, 17).Value = 30
$ws.Cells.Item(8, 4).Value = 44454
$ws.Cells.Item(8, 8).Value = 'Madrigal'
$ws.Cells.Item(8, 10).Value = 100
$ws.Cells.Item(8, 11).Value = 13000
$ws.Cells.Item(8, 12).Value = 14000
$ws.Cells.Item(8, 13).Value = 13500
$ws.Cells.Item(8, 14).Value = '$/caja 40 unidades'
$ws.Cells.Item(8, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(8, 16).Value = 338
$ws.Cells.Item(8, 17).Value = 40
$ws.Cells.Item(9, 8).Value = 'Argentina(o)'
$ws.Cells.Item(9, 10).Value = 180
$ws.Cells.Item(9, 11).Value = 6500
$ws.Cells.Item(9, 12).Value = 7000
$ws.Cells.Item(9, 13).Value = 6778
$ws.Cells.Item(9, 14).Value = '$/caja 50 unidades'
$ws.Cells.Item(9, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(9, 16).Value = 136
$ws.Cells.Item(9, 17).Value = 50
$ws.Cells.Item(10, 4).Value = 44497
$ws.Cells.Item(10, 8).Value = 'Española'
$ws.Cells.Item(10, 10).Value = 200
$ws.Cells.Item(10, 11).Value = 7000
$ws.Cells.Item(10, 12).Value = 7500
$ws.Cells.Item(10, 13).Value = 7250
$ws.Cells.Item(10, 14).Value = '$/caja 30 unidades'
$ws.Cells.Item(10, 16).Value = 242
$ws.Cells.Item(10, 17).Value = 30
$ws.Cells.Item(11, 4).Value = 44497
$ws.Cells.Item(11, 8).Value = 'Madrigal'
$ws.Cells.Item(11, 10).Value = 130
$ws.Cells.Item(11, 11).Value = 6000
$ws.Cells.Item(11, 12).Value = 6500
$ws.Cells.Item(11, 13).Value = 6192
$ws.Cells.Item(11, 14).Value = '$/caja 40 unidades'
$ws.Cells.Item(11, 15).Value = 'Región de Coquimbo'
$ws.Cells.Item(11, 16).Value = 155
$ws.Cells.Item(11, 17).Value = 40
$ws.Cells.Item(12, 4).Value = 44433
$ws.Cells.Item(12, 8).Value = 'Argentina(o)'
$ws.Cells.Item(12, 11).Value = 14000
$ws.Cells.Item(12, 12).Value = 15000
$ws.Cells.Item(12, 13).Value = 14500
$ws.Cells.Item(12, 14).Value = '$/caja 50 unidades'
$ws.Cells.Item(12, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(12, 16).Value = 290
$ws.Cells.Item(12, 17).Value = 50
$ws.Cells.Item(13, 4).Value = 44350
$ws.Cells.Item(13, 8).Value = 'Argentina(o)'
$ws.Cells.Item(13, 10).Value = 50
$ws.Cells.Item(13, 11).Value = 15000
$ws.Cells.Item(13, 12).Value = 16000
$ws.Cells.Item(13, 13).Value = 15600
$ws.Cells.Item(13, 14).Value = '$/caja 50 unidades'
$ws.Cells.Item(13, 16).Value = 312
$ws.Cells.Item(13, 17).Value = 50
$ws.Cells.Item(14, 4).Value = 44350
$ws.Cells.Item(14, 10).Value = 40
$ws.Cells.Item(14, 11).Value = 17000
$ws.Cells.Item(14, 12).Value = 18000
$ws.Cells.Item(14, 13).Value = 17500
$ws.Cells.Item(14, 16).Value = 583
$ws.Cells.Item(15, 4).Value = 44399
$ws.Cells.Item(15, 11).Value = 14000
$ws.Cells.Item(15, 13).Value = 14500
$ws.Cells.Item(15, 16).Value = 483
$ws.Cells.Item(16, 4).Value = 44397
$ws.Cells.Item(16, 8).Value = 'Española'
$ws.Cells.Item(16, 10).Value = 100
$ws.Cells.Item(16, 11).Value = 14000
$ws.Cells.Item(16, 12).Value = 15000
$ws.Cells.Item(16, 13).Value = 14500
$ws.Cells.Item(16, 14).Value = '$/caja 30 unidades'
$ws.Cells.Item(16, 16).Value = 483
$ws.Cells.Item(16, 17).Value = 30
$ws.Cells.Item(17, 8).Value = 'Argentina(o)'
$ws.Cells.Item(17, 11).Value = 8000
$ws.Cells.Item(17, 12).Value = 9000
$ws.Cells.Item(17, 13).Value = 8455
$ws.Cells.Item(17, 14).Value = '$/caja 50 unidades'
$ws.Cells.Item(17, 16).Value = 169
$ws.Cells.Item(17, 17).Value = 50
$ws.Cells.Item(18, 4).Value = 44484
$ws.Cells.Item(18, 8).Value = 'Española'
$ws.Cells.Item(18, 10).Value = 220
$ws.Cells.Item(18, 11).Value = 7500
$ws.Cells.Item(18, 12).Value = 8000
$ws.Cells.Item(18, 13).Value = 7727
$ws.Cells.Item(18, 14).Value = '$/caja 30 unidades'
$ws.Cells.Item(18, 16).Value = 258
$ws.Cells.Item(18, 17).Value = 30
$ws.Cells.Item(19, 4).Value = 44420
$ws.Cells.Item(19, 8).Value = 'Española'
$ws.Cells.Item(19, 14).Value = '$/caja 30 unidades'
$ws.Cells.Item(19, 16).Value = 483
$ws.Cells.Item(19, 17).Value = 30
$ws.Cells.Item(20, 4).Value = 44421
$ws.Cells.Item(20, 10).Value = 100
$ws.Cells.Item(20, 11).Value = 14000
$ws.Cells.Item(20, 12).Value = 15000
$ws.Cells.Item(20, 13).Value = 14500
$ws.Cells.Item(20, 16).Value = 483
$ws.Cells.Item(21, 8).Value = 'Española'
$ws.Cells.Item(21, 10).Value = 450
$ws.Cells.Item(21, 11).Value = 11000
$ws.Cells.Item(21, 12).Value = 12000
$ws.Cells.Item(21, 13).Value = 11444
$ws.Cells.Item(21, 14).Value = '$/caja 30 unidades'
$ws.Cells.Item(21, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(21, 16).Value = 381
$ws.Cells.Item(21, 17).Value = 30
$ws.Cells.Item(22, 4).Value = 44483
$ws.Cells.Item(22, 8).Value = 'Madrigal'
$ws.Cells.Item(22, 10).Value = 220
$ws.Cells.Item(22, 11).Value = 8000
$ws.Cells.Item(22, 12).Value = 8500
$ws.Cells.Item(22, 13).Value = 8273
$ws.Cells.Item(22, 14).Value = '$/caja 40 unidades'
$ws.Cells.Item(22, 15).Value = 'Región de Coquimbo'
$ws.Cells.Item(22, 16).Value = 207
$ws.Cells.Item(22, 17).Value = 40
$ws.Cells.Item(23, 4).Value = 44442
$ws.Cells.Item(23, 11).Value = 14500
$ws.Cells.Item(23, 13).Value = 14750
$ws.Cells.Item(23, 16).Value = 492
$ws.Cells.Item(24, 4).Value = 44463
$ws.Cells.Item(24, 8).Value = 'Argentina(o)'
$ws.Cells.Item(24, 11).Value = 9000
$ws.Cells.Item(24, 12).Value = 10000
$ws.Cells.Item(24, 13).Value = 9500
$ws.Cells.Item(24, 14).Value = '$/caja 50 unidades'
$ws.Cells.Item(24, 16).Value = 190
$ws.Cells.Item(24, 17).Value = 50
$ws.Cells.Item(27, 4).Value = 44447
$ws.Cells.Item(28, 4).Value = 44385
$ws.Cells.Item(30, 4).Value = 44441
$ws.Cells.Item(30, 8).Value = 'Española'
$ws.Cells.Item(30, 10).Value = 100
$ws.Cells.Item(30, 11).Value = 13000
$ws.Cells.Item(30, 12).Value = 14000
$ws.Cells.Item(30, 13).Value = 13500
$ws.Cells.Item(30, 14).Value = '$/caja 30 unidades'
$ws.Cells.Item(30, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(30, 16).Value = 450
$ws.Cells.Item(30, 17).Value = 30
$ws.Cells.Item(31, 4).Value = 44426
$ws.Cells.Item(31, 8).Value = 'Madrigal'
$ws.Cells.Item(31, 11).Value = 12000
$ws.Cells.Item(31, 12).Value = 13000
$ws.Cells.Item(31, 13).Value = 12600
$ws.Cells.Item(31, 14).Value = '$/caja 40 unidades'
$ws.Cells.Item(31, 16).Value = 315
$ws.Cells.Item(31, 17).Value = 40
$ws.Cells.Item(32, 4).Value = 44335
$ws.Cells.Item(32, 8).Value = 'Española'
$ws.Cells.Item(32, 11).Value = 17000
$ws.Cells.Item(32, 12).Value = 18000
$ws.Cells.Item(32, 13).Value = 17500
$ws.Cells.Item(32, 14).Value = '$/caja 30 unidades'
$ws.Cells.Item(32, 16).Value = 583
$ws.Cells.Item(32, 17).Value = 30
$ws.Cells.Item(33, 4).Value = 44435
$ws.Cells.Item(33, 8).Value = 'Argentina(o)'
$ws.Cells.Item(33, 11).Value = 14000
$ws.Cells.Item(33, 12).Value = 15000
$ws.Cells.Item(33, 13).Value = 14500
$ws.Cells.Item(33, 14).Value = '$/caja 50 unidades'
$ws.Cells.Item(33, 16).Value = 290
$ws.Cells.Item(33, 17).Value = 50
$ws.Cells.Item(36, 4).Value = 44383
$ws.Cells.Item(36, 8).Value = 'Argentina(o)'
$ws.Cells.Item(36, 10).Value = 50
$ws.Cells.Item(36, 11).Value = 17000
$ws.Cells.Item(36, 12).Value = 18000
$ws.Cells.Item(36, 13).Value = 17400
$ws.Cells.Item(36, 14).Value = '$/caja 50 unidades'
$ws.Cells.Item(36, 16).Value = 348
$ws.Cells.Item(36, 17).Value = 50
